# Upload new version with timestamp
# Rebuilds the "day sale" shortage report:
#  - rows 7-11 get new item data
#  - 5 brand-new item rows are inserted (old rows 12/13 shift down to 17/18)
#  - the grand-total cell and the generated-at timestamp are refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert five fresh rows after the current last data row (row 11).
#    This pushes the totals row (12) and footer row (13) down to 17/18,
#    exactly matching the target layout.
# ---------------------------------------------------------------------
$ws.Rows("12:16").Insert()

# Give every inserted row the same look (borders/fill/fonts/number
# formats/merges) as the existing data rows, copying per-column from
# row 11 so each cell lands on the same style the template already uses.
foreach ($r in 12..16) {
    $ws.Range("A11:Q11").Copy()
    $ws.Range("A$r`:Q$r").PasteSpecial(-4122)

    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

$excel.CutCopyMode = 0

# Row heights, matching the target XML.
$ws.Rows("12:12").RowHeight = 25.5
$ws.Rows("13:13").RowHeight = 24.75
$ws.Rows("14:14").RowHeight = 25.5
$ws.Rows("15:15").RowHeight = 24.75
$ws.Rows("16:16").RowHeight = 25.5

# ---------------------------------------------------------------------
# 2. Write the item data for rows 7-16.
#    A = row number, C = item name, H = current balance ratio,
#    L = reorder-count, N = price, P = sell price, Q = transaction count.
# ---------------------------------------------------------------------
$items = @(
    @{ Row = 7;  A = 1;  C = "ALKAPRESS PLUS 10/160MG 20 F.C. TABS."; H = "0:1"; L = "1"; N = "102.00"; P = "51.0000";  Q = "0:1" },
    @{ Row = 8;  A = 2;  C = "BISOLOCK-D 5/12.5MG 20 F.C.TAB";        H = "0:0"; L = "1"; N = "28.00";  P = "28.0000";  Q = "1:0" },
    @{ Row = 9;  A = 3;  C = "HALONACE 5 MG 10 TABS.";                H = "0:0"; L = "1"; N = "17.00";  P = "34.0000";  Q = "2:0" },
    @{ Row = 10; A = 4;  C = "MOVXIR 50/500  TAB ";                   H = "1:1"; L = "0"; N = "78.00";  P = "39.0000";  Q = "0:1" },
    @{ Row = 11; A = 5;  C = "NEROTONEX 20 CAPS";                     H = "0:0"; L = "0"; N = "225.00"; P = "112.5000"; Q = "0:1" },
    @{ Row = 12; A = 6;  C = "PANTOMERICAN 40 MG 14 F.C. TABS.";      H = "0:1"; L = "1"; N = "98.00";  P = "49.0000";  Q = "0:1" },
    @{ Row = 13; A = 7;  C = "SPINOBAC 5MG/5ML SYRUP 120ML";          H = "0:0"; L = "1"; N = "39.00";  P = "39.0000";  Q = "1:0" },
    @{ Row = 14; A = 8;  C = "حامل زراع زولا";                        H = "2:0"; L = "0"; N = "40.00";  P = "40.0000";  Q = "1:0" },
    @{ Row = 15; A = 9;  C = "سرنجات 10 سم";                          H = "0:0"; L = "0"; N = "4.00";   P = "4.0000";   Q = "1:0" },
    @{ Row = 16; A = 10; C = "سرنجات 3 سم";                           H = "0:0"; L = "0"; N = "2.00";   P = "2.0000";   Q = "1:0" }
)

# These columns carry numeric-looking number formats (#,##0.##, 0.00, …)
# but the source data is literal text (e.g. "51.0000", "1"), so each value
# is written with a leading apostrophe to force text storage instead of
# letting Excel auto-convert it to a number.
foreach ($item in $items) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("C$r").Value = "'" + $item.C
    $ws.Range("H$r").Value = "'" + $item.H
    $ws.Range("L$r").Value = "'" + $item.L
    $ws.Range("N$r").Value = "'" + $item.N
    $ws.Range("P$r").Value = "'" + $item.P
    $ws.Range("Q$r").Value = "'" + $item.Q
}

# ---------------------------------------------------------------------
# 3. Updated grand total (sum of the P column "sell price" figures).
# ---------------------------------------------------------------------
$ws.Range("P17").Value = 398.5

# ---------------------------------------------------------------------
# 4. Refresh the "generated at" timestamp in the footer.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Monday, 16 June, 2025 10:43 AM"
